$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formulas in E5:E11 (which referenced E2/E3/E4) with new, distinct
# hex colour codes for each ORO type (oro_type rows).
$ws.Range("E5").Value = "#026996"
$ws.Range("E6").Value = "#0688c2"
$ws.Range("E7").Value = "#9ed7f0"
$ws.Range("E8").Value = "#43b08a"
$ws.Range("E9").Value = "#078257"
$ws.Range("E10").Value = "#600787"
$ws.Range("E11").Value = "#ad5ad1"
